$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos price/volume table (GitHub Actions data pull).
# Column D (Price) cells are stored as plain text in the source sheet, even
# though many values look numeric. Purely-numeric replacements are written
# with a leading apostrophe so Excel keeps them as text instead of silently
# converting them to numbers (losing e.g. trailing zeros / formatting).
$ws.Range("D2").Value = '63.242.67'
$ws.Range("E2").Value = '  +0.71%  '
$ws.Range("D3").Value = '2.664.29'
$ws.Range("E3").Value = '  +3.62%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '''605.64'
$ws.Range("E5").Value = '  +4.35%  '
$ws.Range("D6").Value = '''143.25'
$ws.Range("E6").Value = '  -0.25%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").Value = '''0.585'
$ws.Range("E8").Value = '  -0.89%  '
$ws.Range("D9").Value = '2.664.08'
$ws.Range("E9").Value = '  +3.62%  '
$ws.Range("E10").Value = '  +0.12%  '
$ws.Range("D11").Value = '''5.67'
$ws.Range("E11").Value = '  +1.62%  '
$ws.Range("D12").Value = '''0.153'
$ws.Range("E12").Value = '  +0.95%  '
$ws.Range("D13").Value = '''0.355'
$ws.Range("E13").Value = '  +1.98%  '
$ws.Range("D14").Value = '''27.28'
$ws.Range("E14").Value = '  +1.14%  '
$ws.Range("D15").Value = '3.145.94'
$ws.Range("E15").Value = '  +3.77%  '
$ws.Range("D16").Value = '63.132.38'
$ws.Range("E16").Value = '  +0.63%  '
$ws.Range("E17").Value = '  +0.02%  '
$ws.Range("D18").Value = '2.673.48'
$ws.Range("E18").Value = '  +3.69%  '
$ws.Range("D19").Value = '''11.41'
$ws.Range("E19").Value = '  +3.16%  '
$ws.Range("D20").Value = '''338.56'
$ws.Range("E20").Value = '  -0.38%  '
$ws.Range("D21").Value = '''4.38'
$ws.Range("E21").Value = '  +1.29%  '
$ws.Range("D22").Value = '''6.85'
$ws.Range("E22").Value = '  +3.52%  '
$ws.Range("E23").Value = '  +0.10%  '
$ws.Range("D24").Value = '''67.61'
$ws.Range("E24").Value = '  +0.32%  '
$ws.Range("D25").Value = '''1.65'
$ws.Range("E25").Value = '  +3.58%  '
$ws.Range("E26").Value = '  -1.22%  '
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("D28").Value = '''8.50'
$ws.Range("E28").Value = '  +3.45%  '
$ws.Range("E29").Value = '  -0.12%  '
$ws.Range("D30").Value = '''535.91'
$ws.Range("E30").Value = '  +18.24%  '
$ws.Range("E31").Value = '  -1.34%  '
$ws.Range("E32").Value = '  +5.62%  '
$ws.Range("E33").Value = '  +9.59%  '
$ws.Range("D34").Value = '0.0₃0808'
$ws.Range("E34").Value = '  +1.54%  '
$ws.Range("D35").Value = '''173.56'
$ws.Range("E35").Value = '  -1.66%  '
$ws.Range("E36").Value = '  +14.83%  '
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("E38").Value = '  +1.41%  '
$ws.Range("D39").Value = '''19.20'
$ws.Range("E39").Value = '  +1.72%  '
$ws.Range("E40").Value = '  +7.66%  '
$ws.Range("D41").Value = '''174.08'
$ws.Range("E41").Value = '  +11.21%  '
$ws.Range("E43").Value = '  +1.45%  '
$ws.Range("E44").Value = '  +5.60%  '
$ws.Range("D45").Value = '''0.0562'
$ws.Range("E45").Value = '  +5.12%  '
$ws.Range("D46").Value = '''0.633'
$ws.Range("E46").Value = '  +0.05%  '
$ws.Range("E47").Value = '  +0.10%  '
$ws.Range("E48").Value = '  +1.94%  '
$ws.Range("D49").Value = '''18.77'
$ws.Range("E49").Value = '  +4.73%  '
$ws.Range("E50").Value = '  +2.69%  '
$ws.Range("D51").Value = '''11.32'
$ws.Range("E51").Value = '  -1.10%  '
